$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "8.300,01 TL - 199,41 TL"
$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 69,62 TL"
$ws.Range("D14").Value = "3.500 TL - 13.500 TL"
